$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 value (was "pc1", becomes "pc1-pc123")
$ws.Range("A2").Value = "pc1-pc123"

# Delete row 3 entirely (was holding "pc2")
$ws.Range("A3").EntireRow.Delete()

# Update selection to C3 (as seen in the diff)
$ws.Range("C3").Select()
